$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update patient/order data values on Sheet1
$ws.Range("A2").Value = 6379948639
$ws.Range("B2").Value = 600003

# Move the active selection from A2 to I1 (reflects the next UI step in the flow)
$ws.Activate()
$ws.Range("I1").Select()
